$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new shared strings ("line7", "line8") - done implicitly by writing to B8/B9 below.

# Update rows 8-14 (existing rows) - shift labels down: row8/9 become line7/line8,
# rows 10-15 become extr1..extr6 (the old extr1..extr8 block was pushed down by 2
# new rows and relabeled).
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# New rows 16 and 17 (A column uses the same bold/bordered/centered style as
# the other rows' column-A cells, e.g. A2:A15)
$ws.Range("A16").Value = 14
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("A16").VerticalAlignment = -4160
$ws.Range("A16").Borders.LineStyle = 1
$ws.Range("A16").Borders.Weight = 2
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("A17").VerticalAlignment = -4160
$ws.Range("A17").Borders.LineStyle = 1
$ws.Range("A17").Borders.Weight = 2
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
